$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) keeps its original text formatting so
# numeric-looking values like "54.08" or "0.0000330" are not
# coerced into floating point numbers by the COM layer.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '72.693.08'
$ws.Range("E2").Value = '  +1.45%  '
$ws.Range("D3").Value = '4.051.31'
$ws.Range("E3").Value = '  +1.09%  '
$ws.Range("E4").Value = '  -0.16%  '
$ws.Range("D5").Value = '545.02'
$ws.Range("E5").Value = '  +1.95%  '
$ws.Range("D6").Value = '152.70'
$ws.Range("E6").Value = '  +0.87%  '
$ws.Range("D7").Value = '4.043.89'
$ws.Range("E7").Value = '  +1.13%  '
$ws.Range("E8").Value = '  +1.50%  '
$ws.Range("E9").Value = '  -0.04%  '
$ws.Range("D10").Value = '0.757'
$ws.Range("E10").Value = '  +1.23%  '
$ws.Range("E11").Value = '  +0.92%  '
$ws.Range("D12").Value = '54.08'
$ws.Range("E12").Value = '  +13.36%  '
$ws.Range("D13").Value = '0.0000330'
$ws.Range("E13").Value = '  +2.03%  '
$ws.Range("D14").Value = '10.98'
$ws.Range("E14").Value = '  +2.45%  '
$ws.Range("D15").Value = '4.692.49'
$ws.Range("E15").Value = '  +0.83%  '
$ws.Range("D16").Value = '4.043.29'
$ws.Range("E16").Value = '  +1.23%  '
$ws.Range("D17").Value = '14.41'
$ws.Range("E17").Value = '  +2.65%  '
$ws.Range("D18").Value = '20.80'
$ws.Range("E18").Value = '  +1.40%  '
$ws.Range("E19").Value = '  +1.91%  '
$ws.Range("E20").Value = '  -0.50%  '
$ws.Range("D21").Value = '72.522.36'
$ws.Range("E21").Value = '  +1.45%  '
$ws.Range("D22").Value = '451.45'
$ws.Range("E22").Value = '  +5.03%  '
$ws.Range("D23").Value = '98.40'
$ws.Range("E23").Value = '  +0.15%  '
$ws.Range("D24").Value = '3.55'
$ws.Range("E24").Value = '  +0.80%  '
$ws.Range("D25").Value = '4.30'
$ws.Range("E25").Value = '  +2.98%  '
$ws.Range("D26").Value = '14.73'
$ws.Range("E26").Value = '  +2.09%  '
$ws.Range("E27").Value = '  +14.93%  '
$ws.Range("D28").Value = '11.29'
$ws.Range("E28").Value = '  +1.74%  '
$ws.Range("D29").Value = '10.88'
$ws.Range("E29").Value = '  +1.73%  '
$ws.Range("E30").Value = '  +2.05%  '
$ws.Range("D31").Value = '37.39'
$ws.Range("E31").Value = '  +1.91%  '
$ws.Range("D32").Value = '7.99'
$ws.Range("E32").Value = '  +16.51%  '
$ws.Range("D33").Value = '0.135'
$ws.Range("E33").Value = '  +3.99%  '
$ws.Range("D34").Value = '13.66'
$ws.Range("E34").Value = '  +2.28%  '
$ws.Range("D35").Value = '49.07'
$ws.Range("E35").Value = '  +17.63%  '
$ws.Range("D36").Value = '681.10'
$ws.Range("E36").Value = '  -0.29%  '
$ws.Range("D37").Value = '67.14'
$ws.Range("E37").Value = '  +2.41%  '
$ws.Range("D38").Value = '0.451'
$ws.Range("E38").Value = '  +6.62%  '
$ws.Range("D39").Value = '0.0₃0882'
$ws.Range("E39").Value = '  +7.11%  '
$ws.Range("D40").Value = '3.44'
$ws.Range("E40").Value = '  -2.87%  '
$ws.Range("E41").Value = '  -3.24%  '
$ws.Range("B42").Value = 'THORChain'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D42").Value = '11.29'
$ws.Range("E42").Value = '  +18.91%  '
$ws.Range("B43").Value = 'ThetaToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D43").Value = '3.39'
$ws.Range("E43").Value = '  -1.86%  '
$ws.Range("E44").Value = '  +0.10%  '
$ws.Range("E45").Value = '  +2.46%  '
$ws.Range("D46").Value = '0.998'
$ws.Range("E46").Value = '  -0.07%  '
$ws.Range("E47").Value = '  +1.46%  '
$ws.Range("D48").Value = '2.70'
$ws.Range("E48").Value = '  +2.65%  '
$ws.Range("D49").Value = '3.59'
$ws.Range("E49").Value = '  +8.61%  '
$ws.Range("D50").Value = '3.11'
$ws.Range("E50").Value = '  +3.39%  '
$ws.Range("D51").Value = '3.32'
$ws.Range("E51").Value = '  -1.06%  '
